$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / account holder info ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay stored as TEXT (it already
# exceeds Excel's 15 significant digit numeric precision, and the source
# cell type was a string). Force text storage, then restore the original
# cell formatting (style) by pasting formats from a still-unmodified
# neighbouring cell that shares the same style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 03.02.2024"

# --- Row 6 (existing transaction, update values) ---
$ws.Range("B6").Value = "07.02."
$ws.Range("C6").Value = "08.02."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 15905974"
$ws.Range("E6").Value = "86,53-"

# --- Row 7 (existing transaction, update values) ---
$ws.Range("B7").Value = "08.02."
$ws.Range("C7").Value = "09.02."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 78508555"
$ws.Range("E7").Value = "42,25-"

# --- Row 8 (existing transaction, update values) ---
$ws.Range("B8").Value = "09.02."
$ws.Range("C8").Value = "10.02."
$ws.Range("D8").Value = "AMAZON.DE MKTPLC EU SMVDSH"
$ws.Range("E8").Value = "13,75-"

# --- Row 9 (was empty, now a new transaction) ---
# First copy formatting for E9 from E6 so the number style matches (s=17)
$ws.Range("E6").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("B9").Value = "12.02."
$ws.Range("C9").Value = "13.02."
$ws.Range("D9").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E9").Value = "82,32-"

# --- Row 10 (was empty, now a new transaction) ---
$ws.Range("E6").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("B10").Value = "16.02."
$ws.Range("C10").Value = "17.02."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "25,38-"

# --- Row 11 (was empty, now a new transaction) ---
$ws.Range("E6").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("B11").Value = "20.02."
$ws.Range("C11").Value = "21.02."
$ws.Range("D11").Value = "KARTENZ./20.02 LIDL RO"
$ws.Range("E11").Value = "51,52-"

# --- Closing balance date & amount ---
$ws.Range("D12").Value = "KONTOSTAND AM 24.02.2024"
$ws.Range("E12").Value = "301,75-"

# --- Next billing date note ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 05.03.2024"
